$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Test - Test - 10/30/2020" values in column B to "Test - Test"
$ws.Range("B2").Value = "Test - Test"
$ws.Range("B4").Value = "Test - Test"
$ws.Range("B8").Value = "Test - Test"
$ws.Range("B12").Value = "Test - Test"
$ws.Range("B16").Value = "Test - Test"

# Update the selection to C16
$ws.Range("C16").Select()

# Recalculate best-fit width for column B since its content changed
$ws.Columns("B").AutoFit() | Out-Null
